$d = $word.ActiveDocument
$d.Content.Find.Execute("TEAM CAMPUS CONNECT", $true, $false, $false, $false, $false, $true, 1, $false, "TEAM NIGHT SHADES", 2)
